$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = 'se'
$ws.Range("I2").Value = 'Meny'
$ws.Range("I3").Value = 'Nytt spel'
$ws.Range("I4").Value = 'Inställningar'
$ws.Range("C5").Value = 'FPS (Batterieauswirkung):'
$ws.Range("D5").Value = 'FPS (impact de la batterie) :'
$ws.Range("F5").Value = 'FPS (batterij impact):'
$ws.Range("G5").Value = 'FPS (impatto sulla batteria):'
$ws.Range("H5").Value = 'FPS（バッテリーへの影響）:'
$ws.Range("I5").Value = 'FPS (kan påverka batteriet):'
$ws.Range("B6").Value = 'Difficulty (next game):'
$ws.Range("C6").Value = 'Schwierigkeit (Nächstes Spiel):'
$ws.Range("D6").Value = 'Difficulté (prochain jeu) :'
$ws.Range("E6").Value = 'Dificultad (próximo juego):'
$ws.Range("F6").Value = 'Moeilijkheidsgraad (volgende spel):'
$ws.Range("G6").Value = 'Difficoltà (prossima partita):'
$ws.Range("H6").Value = '難易度 (次のゲーム):'
$ws.Range("I6").Value = 'Svårighet (nästa spel):'
$ws.Range("I7").Value = 'Nybörjare'
$ws.Range("I10").Value = 'Språk:'
$ws.Range("I11").Value = 'Engelska'
$ws.Range("I12").Value = 'Tyska'
$ws.Range("I13").Value = 'Franska'
$ws.Range("I14").Value = 'Spanska'
$ws.Range("I15").Value = 'Holländska'
$ws.Range("I16").Value = 'Italienska'
$ws.Range("I17").Value = 'Japanska'
$ws.Range("A18").Value = 'Swedish'
$ws.Range("B18").Value = 'Swedish'
$ws.Range("C18").Value = 'Schwedisch'
$ws.Range("D18").Value = 'Suédois'
$ws.Range("E18").Value = 'Sueco'
$ws.Range("F18").Value = 'Zweeds'
$ws.Range("G18").Value = 'Svedese'
$ws.Range("H18").Value = 'スウェーデン語'
$ws.Range("I18").Value = 'Svenska'
$ws.Range("B21").Value = 'Designed and programmed\nby\nFabrice Bäder'
$ws.Range("A21").Value = 'credits_text'
$ws.Range("C21").Value = 'Entworfen und programmiert\nvon\nFabrice Bäder'
$ws.Range("D21").Value = 'Conçu et programmé\npar\nFabrice bäder'
$ws.Range("E21").Value = 'Diseñado y programado\npor\nFabrice Bäder'
$ws.Range("F21").Value = 'Ontworpen en geprogrammeerd\ndoor\nFabrice Bäder'
$ws.Range("G21").Value = 'Progettato e programmato\nda\nFabrice Bäder'
$ws.Range("I21").Value = 'Designad och programmerad\nav\nFabrice Bäder'
$ws.Range("H21").Value = 'Fabrice Bäder\nによる\nデザインとプログラム'
$ws.Range("I22").Value = 'Fortsätt spelet'
$ws.Range("B23").Value = 'Go'
$ws.Range("C23").Value = 'Los'
$ws.Range("D23").Value = 'Allez'
$ws.Range("E23").Value = 'Ir'
$ws.Range("F23").Value = 'Gaan'
$ws.Range("G23").Value = 'Vai'
$ws.Range("H23").Value = '囲碁'
$ws.Range("I23").Value = 'Gå'
$ws.Range("A24").Value = 'Points: %d'
$ws.Range("B24").Value = 'Points: %d'
$ws.Range("C24").Value = 'Punkte: %d'
$ws.Range("D24").Value = 'Points: %d'
$ws.Range("E24").Value = 'Puntos: %d'
$ws.Range("F24").Value = 'Punten: %d'
$ws.Range("G24").Value = 'Punti: %d'
$ws.Range("H24").Value = 'ポイント: %d'
$ws.Range("I24").Value = 'Poäng: %d'
$ws.Range("A25").Value = 'lost_text'
$ws.Range("B25").Value = 'oh no\nyou''ve lost\nyou got %d points'
$ws.Range("C25").Value = 'Oh Nein\nDu hast verloren\nErreichte Punkte: %d'
$ws.Range("D25").Value = 'oh non\nvous avez perdu\npoints obtenu: %d'
$ws.Range("E25").Value = 'oh no\nhas perdido\npuntos conseguidos: %d'
$ws.Range("F25").Value = 'oh nee\nje hebt verloren\nbehaalde punten: %d'
$ws.Range("G25").Value = 'oh no\nhai perso\npunti ottenuti: %d'
$ws.Range("H25").Value = 'なんてこった\nあなたは負けました\n達成ポイント: %d'
$ws.Range("I25").Value = 'nej då\ndu har förlorat\ndu fick %d poäng'
$ws.Range("C26").Value = 'Neustarten'
$ws.Range("E26").Value = 'reanudar'
$ws.Range("F26").Value = 'Opnieuw opstarten'
$ws.Range("G26").Value = 'riprova'
$ws.Range("I26").Value = 'Försök igen'
$ws.Range("I27").Value = 'Spelet är pausat'

$ws.Range("I21").Select()
